$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 1.074913676625111
$ws.Range("B4").Value = 4.14401568487659
$ws.Range("B5").Value = 3.462833019567579
$ws.Range("B6").Value = 0.6671812968680912
$ws.Range("B7").Value = -6.164758722681718
$ws.Range("B8").Value = 4.439190964013684
$ws.Range("B9").Value = 3.831321260898735
$ws.Range("B10").Value = 0.7051540842417214
$ws.Range("B11").Value = 0.6040035278082057
$ws.Range("B13").Value = 1.166748954083641
$ws.Range("B14").Value = 2.21154834069659
$ws.Range("B15").Value = 3.234237177393018
$ws.Range("B16").Value = 1.234177215189858
$ws.Range("B17").Value = 0.7261884723591683
$ws.Range("B18").Value = -4.719616128339188
$ws.Range("B19").Value = 4.008819402685915
$ws.Range("B20").Value = 2.252360763152805
$ws.Range("B21").Value = -0.313331919805826
$ws.Range("B22").Value = -0.5482818925178212
